# Update cryptos list with latest prices / volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '27.898.93'
$ws.Range('E2').Value = '  +0.79%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.644.37'
$ws.Range('E3').Value = '  +1.01%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.62%  '

# Row 5: BNB
$ws.Range('D5').Value = '''212.00'
$ws.Range('E5').Value = '  -0.12%  '

# Row 6: XRP
$ws.Range('E6').Value = '  +0.38%  '

# Row 7: USDC
$ws.Range('D7').Value = '''0.997'
$ws.Range('E7').Value = '  -0.57%  '

# Row 8: Solana
$ws.Range('D8').Value = '''23.40'
$ws.Range('E8').Value = '  +2.15%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  +1.16%  '

# Row 11: TRON
$ws.Range('D11').Value = '''0.0866'
$ws.Range('E11').Value = '  -2.63%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.872.36'
$ws.Range('E12').Value = '  +0.68%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.640.07'
$ws.Range('E13').Value = '  +0.71%  '

# Row 14: Polkadot
$ws.Range('E14').Value = '  +0.43%  '

# Row 15: Polygon
$ws.Range('D15').Value = '''0.565'
$ws.Range('E15').Value = '  +1.75%  '

# Row 16: Litecoin
$ws.Range('D16').Value = '''65.52'
$ws.Range('E16').Value = '  +1.67%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '27.836.63'
$ws.Range('E17').Value = '  +0.52%  '

# Row 18: BitcoinCash
$ws.Range('D18').Value = '''232.27'
$ws.Range('E18').Value = '  +1.13%  '

# Row 19: Chainlink
$ws.Range('D19').Value = '''7.67'
$ws.Range('E19').Value = '  +1.09%  '

# Row 20: ShibaInu
$ws.Range('D20').Value = '0.0₃0721'

# Row 21: Dai
$ws.Range('D21').Value = '''0.998'
$ws.Range('E21').Value = '  -0.64%  '

# Row 22: Avalanche
$ws.Range('D22').Value = '''10.73'
$ws.Range('E22').Value = '  +7.72%  '

# Row 23: Uniswap
$ws.Range('E23').Value = '  +2.05%  '

# Row 24: Toncoin
$ws.Range('E24').Value = '  +3.17%  '

# Row 25: Monero
$ws.Range('D25').Value = '''150.23'
$ws.Range('E25').Value = '  +0.29%  '

# Row 26: Cosmos
$ws.Range('D26').Value = '''6.92'
$ws.Range('E26').Value = '  +0.39%  '

# Row 27: EthereumClassic
$ws.Range('E27').Value = '  +0.81%  '

# Row 28: Stellar
$ws.Range('E28').Value = '  -0.14%  '

# Row 29: BinanceUSD
$ws.Range('D29').Value = '''0.998'
$ws.Range('E29').Value = '  -0.70%  '

# Row 30: PancakeSwap
$ws.Range('E30').Value = '  +0.11%  '

# Row 31: Hedera
$ws.Range('D31').Value = '''0.0482'
$ws.Range('E31').Value = '  -0.05%  '

# Row 32: Filecoin
$ws.Range('E32').Value = '  +0.16%  '

# Row 33: Maker
$ws.Range('D33').Value = '1.469.94'
$ws.Range('E33').Value = '  +0.51%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('E34').Value = '  +0.31%  '

# Row 35: LidoDAOToken
$ws.Range('E35').Value = '  +1.00%  '

# Row 36: HuobiToken
$ws.Range('E36').Value = '  -1.94%  '

# Row 37: TrustWalletToken
$ws.Range('D37').Value = '''0.938'
$ws.Range('E37').Value = '  +2.74%  '

# Row 38: ARBITRUM
$ws.Range('D38').Value = '''0.888'
$ws.Range('E38').Value = '  +1.97%  '

# Row 39: VeChain
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '''0.560'
$ws.Range('E39').Value = '  -1.29%  '

# Row 40: ImmutableX
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.0168'
$ws.Range('E40').Value = '  +0.27%  '

# Row 41: Aave
$ws.Range('D41').Value = '''69.25'
$ws.Range('E41').Value = '  -0.17%  '

# Row 42: WEMIXToken
$ws.Range('E42').Value = '  -0.26%  '

# Row 43: PaxDollar
$ws.Range('E43').Value = '  -0.67%  '

# Row 44: MXToken
$ws.Range('B44').Value = 'mCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D44').Value = '''2.45'
$ws.Range('E44').Value = '  -0.47%  '

# Row 45: mCoin
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '''2.22'
$ws.Range('E45').Value = '  -0.28%  '

# Row 46: FraxShare
$ws.Range('D46').Value = '''5.38'
$ws.Range('E46').Value = '  -0.95%  '

# Row 47: RocketPoolETH
$ws.Range('D47').Value = '1.785.30'
$ws.Range('E47').Value = '  +0.58%  '

# Row 48: RenderToken
$ws.Range('E48').Value = '  +3.75%  '

# Row 49: Quant
$ws.Range('D49').Value = '''87.83'
$ws.Range('E49').Value = '  +2.34%  '

# Row 50: Algorand
$ws.Range('E50').Value = '  +1.60%  '

# Row 51: BabyDogeCoin
$ws.Range('D51').Value = '0.0₆0101'
$ws.Range('E51').Value = '  +1.89%  '
